$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 386 (shifts old 386..395 down to 388..397)
$ws.Range("A386:T387").EntireRow.Insert()

# New row 386 (date 2023-08-09 / 45147, "Primera")
$ws.Range("A386").Value = 6
$ws.Range("B386").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C386").Value = "Metropolitana"
$ws.Range("D386").Value = 45147
$ws.Range("E386").Value = 13
$ws.Range("F386").Value = "Fruta"
$ws.Range("G386").Value = 100107
$ws.Range("H386").Value = "Otros"
$ws.Range("I386").Value = 100107002
$ws.Range("J386").Value = "Chirimoya"
$ws.Range("K386").Value = "Cultivar IV Región"
$ws.Range("L386").Value = "Primera"
$ws.Range("M386").Value = 175
$ws.Range("N386").Value = 24000
$ws.Range("O386").Value = 24000
$ws.Range("P386").Value = 24000
$ws.Range("Q386").Value = "`$/bandeja 10 kilos"
$ws.Range("R386").Value = "Provincia de Limarí"
$ws.Range("S386").Value = 2400
$ws.Range("T386").Value = 10

# New row 387 (date 2023-08-09 / 45147, "Segunda")
$ws.Range("A387").Value = 6
$ws.Range("B387").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C387").Value = "Metropolitana"
$ws.Range("D387").Value = 45147
$ws.Range("E387").Value = 13
$ws.Range("F387").Value = "Fruta"
$ws.Range("G387").Value = 100107
$ws.Range("H387").Value = "Otros"
$ws.Range("I387").Value = 100107002
$ws.Range("J387").Value = "Chirimoya"
$ws.Range("K387").Value = "Cultivar IV Región"
$ws.Range("L387").Value = "Segunda"
$ws.Range("M387").Value = 150
$ws.Range("N387").Value = 21000
$ws.Range("O387").Value = 21000
$ws.Range("P387").Value = 21000
$ws.Range("Q387").Value = "`$/bandeja 10 kilos"
$ws.Range("R387").Value = "Provincia de Limarí"
$ws.Range("S387").Value = 2100
$ws.Range("T387").Value = 10
